# "Generate Report for Handoff"
# Updates the localization-status report to reflect that the files are now
# "Ready for handoff" (was "In Translation") and refreshes the handoff
# timestamps, then re-sizes the Status columns to fit the new, longer text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Update status text: "In Translation" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Update the handoff generation timestamps ---
$wsZhCn.Range("H2").Value = "2016-08-18 09:00:36"
$wsDeDe.Range("H2").Value = "2016-08-18 09:00:41"
$wsOverview.Range("G2").Value = "2016-08-18 09:00:41"

# --- Resize the "Status" column(s) to fit the new, longer text ---
# (target characters-width is ~17.216; the host rounds ColumnWidth to the
# nearest whole pixel, so we pick the input that lands closest to it)
$wsOverview.Columns.Item(5).ColumnWidth = 16.3333333
$wsOverview.Columns.Item(6).ColumnWidth = 16.3333333
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3333333
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3333333
